# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps advance
# - The zh-cn / de-de "Status" + "Latest Handoff Datetime" columns widen to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet -----------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 06:46:14"

# Widen the zh-cn / de-de status columns on the Overview sheet so the
# longer "Ready for handoff" text fits.
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 16.33
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 16.33

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 06:46:08"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 16.33

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 06:46:14"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 16.33
